$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (Word's "last edit location" marker).
#    It currently sits right after the run containing "6" (inside the "...de {s1f6}" text).
#    It will be re-created below at the new last-edit location once the text edit is made.
$d.Bookmarks("_GoBack").Delete()

# 2. Replace the literal placeholder text with the {o1} merge field that drives the
#    "original checked/unchecked" logic.
$findRange = $d.Content
$findRange.Find.Execute("Copia Verdadera del Original", $true, $false, $false, $false, $false, $true, 1, $false, "{o1}", 2)

# 3. Re-create the "_GoBack" bookmark collapsed immediately after the new "{o1}" run, matching
#    where Word leaves the edit-location bookmark after typing/replacing text there.
#    A temporary marker character is used to work around the fact that collapsing a range to a
#    position that is simultaneously "end of run" and "end of paragraph" can otherwise resolve
#    to the wrong place; the marker is removed again once the bookmark is anchored.
$insPoint = $d.Range($findRange.End, $findRange.End)
$insPoint.InsertAfter("#")

$bmPoint = $d.Range($findRange.End, $findRange.End)
$d.Bookmarks.Add("_GoBack", $bmPoint)

$marker = $d.Range($findRange.End, $findRange.End + 1)
$marker.Delete()
